$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-47 (reshuffled data values) ---
$ws.Cells.Item(2, 4).Value = 44400
$ws.Cells.Item(2, 13).Value = 5
$ws.Cells.Item(2, 14).Value = 24000
$ws.Cells.Item(2, 15).Value = 24000
$ws.Cells.Item(2, 16).Value = 24000
$ws.Cells.Item(2, 19).Value = 1200
$ws.Cells.Item(3, 4).Value = 44396
$ws.Cells.Item(3, 13).Value = 12
$ws.Cells.Item(4, 4).Value = 44249
$ws.Cells.Item(4, 14).Value = 25000
$ws.Cells.Item(4, 15).Value = 25000
$ws.Cells.Item(4, 16).Value = 25000
$ws.Cells.Item(4, 19).Value = 1250
$ws.Cells.Item(5, 4).Value = 44221
$ws.Cells.Item(5, 13).Value = 30
$ws.Cells.Item(6, 4).Value = 44251
$ws.Cells.Item(6, 13).Value = 15
$ws.Cells.Item(6, 14).Value = 25000
$ws.Cells.Item(6, 15).Value = 25000
$ws.Cells.Item(6, 16).Value = 25000
$ws.Cells.Item(6, 19).Value = 1250
$ws.Cells.Item(7, 4).Value = 44442
$ws.Cells.Item(7, 13).Value = 25
$ws.Cells.Item(7, 14).Value = 23000
$ws.Cells.Item(7, 15).Value = 23000
$ws.Cells.Item(7, 16).Value = 23000
$ws.Cells.Item(7, 19).Value = 1150
$ws.Cells.Item(8, 4).Value = 44398
$ws.Cells.Item(8, 13).Value = 15
$ws.Cells.Item(8, 14).Value = 25000
$ws.Cells.Item(8, 15).Value = 25000
$ws.Cells.Item(8, 16).Value = 25000
$ws.Cells.Item(8, 19).Value = 1250
$ws.Cells.Item(9, 4).Value = 44420
$ws.Cells.Item(9, 13).Value = 35
$ws.Cells.Item(9, 14).Value = 25000
$ws.Cells.Item(9, 15).Value = 25000
$ws.Cells.Item(9, 16).Value = 25000
$ws.Cells.Item(9, 19).Value = 1250
$ws.Cells.Item(10, 4).Value = 44462
$ws.Cells.Item(10, 13).Value = 10
$ws.Cells.Item(10, 14).Value = 24000
$ws.Cells.Item(10, 15).Value = 24000
$ws.Cells.Item(10, 16).Value = 24000
$ws.Cells.Item(10, 19).Value = 1200
$ws.Cells.Item(11, 4).Value = 44431
$ws.Cells.Item(11, 13).Value = 40
$ws.Cells.Item(12, 4).Value = 44238
$ws.Cells.Item(12, 13).Value = 30
$ws.Cells.Item(12, 14).Value = 25000
$ws.Cells.Item(12, 15).Value = 25000
$ws.Cells.Item(12, 16).Value = 25000
$ws.Cells.Item(12, 19).Value = 1250
$ws.Cells.Item(13, 4).Value = 44222
$ws.Cells.Item(13, 13).Value = 15
$ws.Cells.Item(14, 4).Value = 44175
$ws.Cells.Item(14, 13).Value = 25
$ws.Cells.Item(14, 14).Value = 23000
$ws.Cells.Item(14, 15).Value = 23000
$ws.Cells.Item(14, 16).Value = 23000
$ws.Cells.Item(14, 19).Value = 1150
$ws.Cells.Item(15, 4).Value = 44419
$ws.Cells.Item(15, 13).Value = 40
$ws.Cells.Item(16, 4).Value = 44418
$ws.Cells.Item(16, 13).Value = 20
$ws.Cells.Item(16, 14).Value = 24000
$ws.Cells.Item(16, 15).Value = 24000
$ws.Cells.Item(16, 16).Value = 24000
$ws.Cells.Item(16, 19).Value = 1200
$ws.Cells.Item(17, 4).Value = 44424
$ws.Cells.Item(17, 13).Value = 25
$ws.Cells.Item(18, 4).Value = 44412
$ws.Cells.Item(18, 13).Value = 20
$ws.Cells.Item(19, 4).Value = 44428
$ws.Cells.Item(19, 13).Value = 15
$ws.Cells.Item(20, 4).Value = 44421
$ws.Cells.Item(20, 13).Value = 20
$ws.Cells.Item(21, 4).Value = 44232
$ws.Cells.Item(21, 13).Value = 15
$ws.Cells.Item(22, 4).Value = 44356
$ws.Cells.Item(22, 13).Value = 15
$ws.Cells.Item(22, 14).Value = 24000
$ws.Cells.Item(22, 15).Value = 24000
$ws.Cells.Item(22, 16).Value = 24000
$ws.Cells.Item(22, 19).Value = 1200
$ws.Cells.Item(23, 4).Value = 44454
$ws.Cells.Item(23, 13).Value = 25
$ws.Cells.Item(23, 14).Value = 25000
$ws.Cells.Item(23, 15).Value = 25000
$ws.Cells.Item(23, 16).Value = 25000
$ws.Cells.Item(23, 19).Value = 1250
$ws.Cells.Item(24, 4).Value = 44467
$ws.Cells.Item(25, 4).Value = 44235
$ws.Cells.Item(25, 13).Value = 15
$ws.Cells.Item(27, 4).Value = 44469
$ws.Cells.Item(27, 13).Value = 40
$ws.Cells.Item(28, 4).Value = 44434
$ws.Cells.Item(28, 13).Value = 20
$ws.Cells.Item(29, 4).Value = 44214
$ws.Cells.Item(29, 13).Value = 15
$ws.Cells.Item(29, 14).Value = 25000
$ws.Cells.Item(29, 15).Value = 25000
$ws.Cells.Item(29, 16).Value = 25000
$ws.Cells.Item(29, 19).Value = 1250
$ws.Cells.Item(30, 4).Value = 44391
$ws.Cells.Item(30, 13).Value = 10
$ws.Cells.Item(30, 14).Value = 24000
$ws.Cells.Item(30, 15).Value = 24000
$ws.Cells.Item(30, 16).Value = 24000
$ws.Cells.Item(30, 19).Value = 1200
$ws.Cells.Item(31, 4).Value = 44231
$ws.Cells.Item(31, 13).Value = 15
$ws.Cells.Item(32, 4).Value = 44475
$ws.Cells.Item(33, 4).Value = 44474
$ws.Cells.Item(33, 13).Value = 20
$ws.Cells.Item(33, 14).Value = 24000
$ws.Cells.Item(33, 15).Value = 24000
$ws.Cells.Item(33, 16).Value = 24000
$ws.Cells.Item(33, 19).Value = 1200
$ws.Cells.Item(34, 4).Value = 44349
$ws.Cells.Item(34, 13).Value = 30
$ws.Cells.Item(35, 4).Value = 44452
$ws.Cells.Item(35, 13).Value = 25
$ws.Cells.Item(35, 14).Value = 25000
$ws.Cells.Item(35, 15).Value = 25000
$ws.Cells.Item(35, 16).Value = 25000
$ws.Cells.Item(35, 19).Value = 1250
$ws.Cells.Item(36, 4).Value = 44435
$ws.Cells.Item(36, 13).Value = 100
$ws.Cells.Item(36, 14).Value = 24000
$ws.Cells.Item(36, 15).Value = 24000
$ws.Cells.Item(36, 16).Value = 24000
$ws.Cells.Item(36, 19).Value = 1200
$ws.Cells.Item(37, 4).Value = 44433
$ws.Cells.Item(37, 13).Value = 10
$ws.Cells.Item(38, 4).Value = 44334
$ws.Cells.Item(38, 13).Value = 20
$ws.Cells.Item(38, 14).Value = 25000
$ws.Cells.Item(38, 15).Value = 25000
$ws.Cells.Item(38, 16).Value = 25000
$ws.Cells.Item(38, 19).Value = 1250
$ws.Cells.Item(39, 4).Value = 44468
$ws.Cells.Item(39, 13).Value = 20
$ws.Cells.Item(39, 14).Value = 24000
$ws.Cells.Item(39, 15).Value = 24000
$ws.Cells.Item(39, 16).Value = 24000
$ws.Cells.Item(39, 19).Value = 1200
$ws.Cells.Item(40, 4).Value = 44363
$ws.Cells.Item(40, 13).Value = 30
$ws.Cells.Item(41, 4).Value = 44414
$ws.Cells.Item(41, 13).Value = 15
$ws.Cells.Item(42, 4).Value = 44392
$ws.Cells.Item(42, 13).Value = 10
$ws.Cells.Item(43, 4).Value = 44425
$ws.Cells.Item(43, 13).Value = 15
$ws.Cells.Item(44, 4).Value = 44390
$ws.Cells.Item(44, 13).Value = 10
$ws.Cells.Item(45, 4).Value = 44466
$ws.Cells.Item(45, 13).Value = 70
$ws.Cells.Item(45, 14).Value = 24000
$ws.Cells.Item(45, 15).Value = 24000
$ws.Cells.Item(45, 16).Value = 24000
$ws.Cells.Item(45, 19).Value = 1200
$ws.Cells.Item(47, 4).Value = 44389
$ws.Cells.Item(47, 13).Value = 20

# --- Add new row 48 (copy of original row-3 record) ---
$ws.Cells.Item(48, 1).Value = 10
$ws.Cells.Item(48, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(48, 3).Value = "La Araucanía"
$ws.Cells.Item(48, 4).Value = 44432
$ws.Cells.Item(48, 4).NumberFormat = $ws.Cells.Item(47, 4).NumberFormat()
$ws.Cells.Item(48, 5).Value = 9
$ws.Cells.Item(48, 6).Value = "Fruta"
$ws.Cells.Item(48, 7).Value = 100108
$ws.Cells.Item(48, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(48, 9).Value = 100108007
$ws.Cells.Item(48, 10).Value = "Coco"
$ws.Cells.Item(48, 11).Value = "Sin especificar"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value = 30
$ws.Cells.Item(48, 14).Value = 24000
$ws.Cells.Item(48, 15).Value = 24000
$ws.Cells.Item(48, 16).Value = 24000
$ws.Cells.Item(48, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(48, 18).Value = "Perú"
$ws.Cells.Item(48, 19).Value = 1200
$ws.Cells.Item(48, 20).Value = 20
